# "Added Mapping column to import excel"
#
# The import-mapping header row is relabelled (Polish field names used by
# the mapping UI instead of the raw C# property names), and the currently
# selected cell on the sheet becomes E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Re-label the header row (row 1) to the new mapping column names.
# B1 ("SecondName") and E1 ("AddressId") are left unchanged.
$ws.Range("A1").Value = "imie"
$ws.Range("C1").Value = "nazwisko"
$ws.Range("D1").Value = "PESEL"
$ws.Range("F1").Value = "Telefon 1"
$ws.Range("G1").Value = "Telefon2"

# Move the active selection to E6, matching the saved view state.
$ws.Range("E6").Select()

# Best-effort: restore the (cosmetic) saved workbook window geometry.
$win = $excel.ActiveWindow
$win.Left = -15
$win.Top = -15
$win.Width = 21630
$win.Height = 5055
